$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 119, shifting existing rows 119:184 down to 120:185.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with its data.
$ws.Cells.Item(119, 1).Value = 9
$ws.Cells.Item(119, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(119, 3).Value = "Metropolitana"
$ws.Cells.Item(119, 4).Value = 44438
$ws.Cells.Item(119, 5).Value = 13
$ws.Cells.Item(119, 6).Value = 100112032
$ws.Cells.Item(119, 7).Value = "Zapallo italiano"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 106
$ws.Cells.Item(119, 11).Value = 13000
$ws.Cells.Item(119, 12).Value = 14000
$ws.Cells.Item(119, 13).Value = 13500
$ws.Cells.Item(119, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(119, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(119, 16).Value = 225
$ws.Cells.Item(119, 17).Value = 60
$ws.Cells.Item(119, 18).Value = "Hortaliza"
